# Auto-generated Excel COM-interop script applying the Siren_Profits.xlsx diff.
# Updates market-price / profit columns (H-N) across the ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC!row 4
$ws.Range("H4").Value = 1475.25
$ws.Range("I4").Value = 1543.1428
$ws.Range("J4").Value = 1000
$ws.Range("K4").Value = 1543.1428
$ws.Range("L4").Value = 1000
$ws.Range("M4").Value = -1429.1428
$ws.Range("N4").Value = -1228

# ALC!row 69
$ws.Range("H69").Value = 9250
$ws.Range("I69").Value = 9000
$ws.Range("K69").Value = 27000
$ws.Range("M69").Value = -26126

# ALC!row 72
$ws.Range("H72").Value = 9250
$ws.Range("I72").Value = 9000
$ws.Range("K72").Value = 81000
$ws.Range("M72").Value = -76632

# ALC!row 98
$ws.Range("H98").Value = 25930.5
$ws.Range("I98").Value = 43981.77
$ws.Range("K98").Value = 43981.77
$ws.Range("M98").Value = -42483.77

# ALC!row 122
$ws.Range("H122").Value = 25930.5
$ws.Range("I122").Value = 43981.77
$ws.Range("K122").Value = 131945.31
$ws.Range("M122").Value = -129495.31

# ALC!row 125
$ws.Range("H125").Value = 6133.375
$ws.Range("I125").Value = 9169.25
$ws.Range("K125").Value = 82523.25
$ws.Range("M125").Value = -80063.25

# ALC!row 129
$ws.Range("H129").Value = 1610.9333
$ws.Range("I129").Value = 1305.3334
$ws.Range("J129").Value = 2833.3333
$ws.Range("K129").Value = 3916.0002
$ws.Range("L129").Value = 8499.999899999999
$ws.Range("M129").Value = 1083.9998
$ws.Range("N129").Value = -18499.9999

# ALC!row 131
$ws.Range("H131").Value = 2191.7144
$ws.Range("I131").Value = 1710.4
$ws.Range("K131").Value = 5131.200000000001
$ws.Range("M131").Value = -91.20000000000073

# ALC!row 135
$ws.Range("H135").Value = 7083.696
$ws.Range("I135").Value = 8634.294
$ws.Range("K135").Value = 77708.64599999999
$ws.Range("M135").Value = -75173.64599999999

$ws = $wb.Worksheets.Item("ARM")
# ARM!row 32
$ws.Range("H32").Value = 4904.875
$ws.Range("I32").Value = 4468.0156
$ws.Range("K32").Value = 4468.0156
$ws.Range("M32").Value = -4181.0156

# ARM!row 110
$ws.Range("H110").Value = 1248.4445
$ws.Range("I110").Value = 1010.2083
$ws.Range("J110").Value = 3154.3333
$ws.Range("K110").Value = 1010.2083
$ws.Range("L110").Value = 3154.3333
$ws.Range("M110").Value = 1034.7917
$ws.Range("N110").Value = -7244.3333

# ARM!row 122
$ws.Range("H122").Value = 1672154
$ws.Range("I122").Value = 4912.4287
$ws.Range("K122").Value = 14737.2861
$ws.Range("M122").Value = -12287.2861

# ARM!row 133
$ws.Range("H133").Value = 62297.6
$ws.Range("J133").Value = 70312.75
$ws.Range("L133").Value = 70312.75
$ws.Range("N133").Value = -75372.75

# ARM!row 138
$ws.Range("H138").Value = 62027.5
$ws.Range("J138").Value = 62027.5
$ws.Range("L138").Value = 62027.5
$ws.Range("N138").Value = -72307.5

$ws = $wb.Worksheets.Item("BSM")
# BSM!row 86
$ws.Range("H86").Value = 7534.0625
$ws.Range("I86").Value = 11936.875
$ws.Range("J86").Value = 3131.25
$ws.Range("K86").Value = 11936.875
$ws.Range("L86").Value = 3131.25
$ws.Range("M86").Value = -10813.875
$ws.Range("N86").Value = -5377.25

# BSM!row 89
$ws.Range("H89").Value = 7534.0625
$ws.Range("I89").Value = 11936.875
$ws.Range("J89").Value = 3131.25
$ws.Range("K89").Value = 59684.375
$ws.Range("L89").Value = 15656.25
$ws.Range("M89").Value = -54068.375
$ws.Range("N89").Value = -26888.25

# BSM!row 94
$ws.Range("H94").Value = 4180.364
$ws.Range("I94").Value = 3173.077
$ws.Range("K94").Value = 3173.077
$ws.Range("M94").Value = -2722.077

# BSM!row 105
$ws.Range("H105").Value = 69258.766
$ws.Range("I105").Value = 94691.586
$ws.Range("K105").Value = 94691.586
$ws.Range("M105").Value = -92944.586

# BSM!row 138
$ws.Range("H138").Value = 124999.25
$ws.Range("I138").Value = 50000
$ws.Range("J138").Value = 149999
$ws.Range("K138").Value = 50000
$ws.Range("L138").Value = 149999
$ws.Range("N138").Value = -160279
$ws.Range("M138").Value = -44860

$ws = $wb.Worksheets.Item("CRP")
# CRP!row 103
$ws.Range("H103").Value = 14076.071
$ws.Range("I103").Value = 19341
$ws.Range("J103").Value = 4599.2
$ws.Range("K103").Value = 19341
$ws.Range("L103").Value = 4599.2
$ws.Range("M103").Value = -18169
$ws.Range("N103").Value = -6943.2

# CRP!row 107
$ws.Range("H107").Value = 12373.909
$ws.Range("I107").Value = 15951.625
$ws.Range("J107").Value = 2833.3333
$ws.Range("K107").Value = 15951.625
$ws.Range("L107").Value = 2833.3333
$ws.Range("M107").Value = -14031.625
$ws.Range("N107").Value = -6673.3333

$ws = $wb.Worksheets.Item("CUL")
# CUL!row 25
$ws.Range("H25").Value = 2184.7058
$ws.Range("J25").Value = 2468
$ws.Range("L25").Value = 7404
$ws.Range("N25").Value = -7742

# CUL!row 30
$ws.Range("H30").Value = 2184.7058
$ws.Range("J30").Value = 2468
$ws.Range("L30").Value = 7404
$ws.Range("N30").Value = -7608

# CUL!row 39
$ws.Range("H39").Value = 369.44446
$ws.Range("I39").Value = 309.375
$ws.Range("K39").Value = 928.125
$ws.Range("M39").Value = -634.125

# CUL!row 57
$ws.Range("H57").Value = 15000
$ws.Range("I57").Value = 10000
$ws.Range("K57").Value = 30000
$ws.Range("M57").Value = -29441

# CUL!row 62
$ws.Range("H62").Value = 1166.1666
$ws.Range("I62").Value = 999.3333
$ws.Range("J62").Value = 1333
$ws.Range("K62").Value = 2997.9999
$ws.Range("L62").Value = 3999
$ws.Range("M62").Value = -2311.9999
$ws.Range("N62").Value = -5371

# CUL!row 65
$ws.Range("H65").Value = 1166.1666
$ws.Range("I65").Value = 999.3333
$ws.Range("J65").Value = 1333
$ws.Range("K65").Value = 8993.9997
$ws.Range("L65").Value = 11997
$ws.Range("M65").Value = -5561.9997
$ws.Range("N65").Value = -18861

# CUL!row 93
$ws.Range("H93").Value = 24399.8

# CUL!row 94
$ws.Range("H94").Value = 0
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 0
$ws.Range("N94").Value = 0
$ws.Range("L94").ClearContents()
$ws.Range("M94").ClearContents()

# CUL!row 104
$ws.Range("H104").Value = 4525
$ws.Range("I104").Value = 3321.6667
$ws.Range("K104").Value = 9965.000100000001
$ws.Range("M104").Value = -7344.000100000001

# CUL!row 113
$ws.Range("H113").Value = 1752.5
$ws.Range("I113").Value = 1300
$ws.Range("K113").Value = 3900
$ws.Range("M113").Value = -1730

# CUL!row 121
$ws.Range("H121").Value = 4175.8
$ws.Range("J121").Value = 3783
$ws.Range("L121").Value = 11349
$ws.Range("N121").Value = -13969

$ws = $wb.Worksheets.Item("GSM")
# GSM!row 80
$ws.Range("H80").Value = 7035.4707
$ws.Range("I80").Value = 7557.5454
$ws.Range("K80").Value = 7557.5454
$ws.Range("M80").Value = -6559.5454

# GSM!row 83
$ws.Range("H83").Value = 7035.4707
$ws.Range("I83").Value = 7557.5454
$ws.Range("K83").Value = 37787.727
$ws.Range("M83").Value = -32795.727

# GSM!row 122
$ws.Range("H122").Value = 8056.3706
$ws.Range("I122").Value = 5260.88
$ws.Range("K122").Value = 15782.64
$ws.Range("M122").Value = -13332.64

# GSM!row 126
$ws.Range("H126").Value = 6247.3237
$ws.Range("I126").Value = 7883.952
$ws.Range("K126").Value = 23651.856
$ws.Range("M126").Value = -21181.856

# GSM!row 132
$ws.Range("H132").Value = 10930.4
$ws.Range("I132").Value = 13615.714
$ws.Range("K132").Value = 40847.142
$ws.Range("M132").Value = -38317.142

# GSM!row 138
$ws.Range("H138").Value = 96999
$ws.Range("J138").Value = 96999
$ws.Range("L138").Value = 96999
$ws.Range("N138").Value = -107279

# GSM!row 141
$ws.Range("H141").Value = 84332.664
$ws.Range("J141").Value = 84332.664
$ws.Range("L141").Value = 84332.664
$ws.Range("N141").Value = -94692.664

$ws = $wb.Worksheets.Item("LTW")
# LTW!row 82
$ws.Range("H82").Value = 2549.6
$ws.Range("I82").Value = 3283
$ws.Range("J82").Value = 1449.5
$ws.Range("K82").Value = 3283
$ws.Range("L82").Value = 1449.5
$ws.Range("M82").Value = -2922
$ws.Range("N82").Value = -2171.5

# LTW!row 85
$ws.Range("H85").Value = 2549.6
$ws.Range("I85").Value = 3283
$ws.Range("J85").Value = 1449.5
$ws.Range("K85").Value = 3283
$ws.Range("L85").Value = 1449.5
$ws.Range("M85").Value = -2035
$ws.Range("N85").Value = -3945.5

# LTW!row 140
$ws.Range("H140").Value = 120000
$ws.Range("J140").Value = 120000
$ws.Range("L140").Value = 120000
$ws.Range("N140").Value = -130360

$ws = $wb.Worksheets.Item("WVR")
# WVR!row 58
$ws.Range("H58").Value = 8510622
$ws.Range("I58").Value = 14489
$ws.Range("J58").Value = 11342667
$ws.Range("K58").Value = 14489
$ws.Range("L58").Value = 11342667
$ws.Range("M58").Value = -14181
$ws.Range("N58").Value = -11343283

# WVR!row 62
$ws.Range("H62").Value = 385072.56
$ws.Range("I62").Value = 683400.8
$ws.Range("J62").Value = 12162.25
$ws.Range("K62").Value = 683400.8
$ws.Range("L62").Value = 12162.25
$ws.Range("M62").Value = -682776.8
$ws.Range("N62").Value = -13410.25

# WVR!row 65
$ws.Range("H65").Value = 385072.56
$ws.Range("I65").Value = 683400.8
$ws.Range("J65").Value = 12162.25
$ws.Range("K65").Value = 3417004
$ws.Range("L65").Value = 60811.25
$ws.Range("M65").Value = -3413884
$ws.Range("N65").Value = -67051.25

# WVR!row 132
$ws.Range("H132").Value = 3233.535
$ws.Range("I132").Value = 3109.5151
$ws.Range("K132").Value = 9328.5453
$ws.Range("M132").Value = -6798.5453

# WVR!row 139
$ws.Range("H139").Value = 115333.336
$ws.Range("J139").Value = 68000
$ws.Range("L139").Value = 68000
$ws.Range("N139").Value = -78280

